# Case_0_56 / res_bus / vm_pu.xlsx - update bus voltage magnitude results
# (case with 380 kV done): slack bus voltage setpoint changed from 1.05 to
# 1.02 p.u. and all downstream per-unit voltage results recalculated.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array holds the new values for columns B,C,D,E,F,I,J,K,L,M,N
# (columns G and H are untouched) for data rows 2..25.
$data = @(
    ,@(1.02,1.031528711567798,1.040777922395563,1.031004210401216,1.047580054212718,1.02359499962809,1.036663938705474,1.043559307405298,1.033813583220027,1.050342295136892,1.038136120405452)
    ,@(1.02,1.033149473229332,1.042373664865802,1.032405029692096,1.049331019988976,1.023504579208684,1.037923198507259,1.04496347961889,1.03502128702245,1.051902701189664,1.039397168500594)
    ,@(1.02,1.034196746797847,1.04340511938952,1.033310310487999,1.050463109509821,1.023443335729026,1.038736259236988,1.045870500420916,1.035801096219006,1.0529110404373,1.040211383869819)
    ,@(1.02,1.034636679974207,1.043838489932026,1.033690625064755,1.050938834856217,1.023416932628374,1.039077655643758,1.046251444461984,1.036128539092527,1.05333463351988,1.040553265098645)
    ,@(1.02,1.034710527025688,1.043911240204761,1.033754466224418,1.051018699525775,1.023412460900047,1.039134953541415,1.0463153853514,1.036183495603341,1.053405738610759,1.040610644365889)
    ,@(1.02,1.034202626530109,1.043410911088182,1.03331539330175,1.050469466965877,1.02344298551034,1.038740822609223,1.045875592045448,1.035805473044411,1.052916701723883,1.040215953722566)
    ,@(1.02,1.03207676590428,1.041317444217922,1.031477865419341,1.048171994757702,1.023565007130014,1.03708988171624,1.044034184986659,1.034222079782421,1.050869926735621,1.038562668304149)
    ,@(1.02,1.028319006331301,1.037619640883564,1.02823078613257,1.044116151614281,1.023759139049229,1.034166833571057,1.040776944220929,1.031418918133991,1.047252478801198,1.035635469095322)
    ,@(1.02,1.025805277719833,1.035147871614437,1.026059419911654,1.041406588782685,1.023874603872166,1.032208291231637,1.038596517508379,1.029540915074887,1.044832980821995,1.033674145401131)
    ,@(1.02,1.024714625673881,1.034075872470286,1.025117502779489,1.04023182003385,1.023921309042974,1.031357771619474,1.037650123774023,1.028725424021697,1.043783306098854,1.032822417953563)
    ,@(1.02,1.024309166753347,1.033677415882243,1.02476736724914,1.039795218422661,1.023938164268129,1.031041471582674,1.037298241350387,1.028422159079493,1.043393095180531,1.032505668734431)
    ,@(1.02,1.024396154763193,1.033762898510582,1.02484248468764,1.039888882034356,1.023934571058996,1.031109336333286,1.037373737323711,1.028487226663756,1.043476811196233,1.032573629860773)
    ,@(1.02,1.024681117369311,1.034042941468732,1.025088565942884,1.040195735385687,1.023922712353276,1.03133163395268,1.037621044259477,1.028700363347615,1.043751057638611,1.032796243168284)
    ,@(1.02,1.02485664650103,1.034215449256043,1.025240149235507,1.040384765635166,1.023915340507297,1.031468548309034,1.037773371563874,1.028831636541383,1.043919987820462,1.032933351958729)
    ,@(1.02,1.025877613434061,1.035218979689825,1.026121895156959,1.041484521364971,1.023871434968932,1.032264684788369,1.038659278111698,1.029594987217836,1.044902600769328,1.033730619043183)
    ,@(1.02,1.026517443274939,1.035848002466259,1.026674528593646,1.04217395540717,1.023843013862359,1.032763415155289,1.039214372583784,1.030073193787924,1.045518419634303,1.034230057664429)
    ,@(1.02,1.026890434684554,1.036214736947008,1.026996706862236,1.042575945815628,1.023826118601784,1.033054079667692,1.039537933006811,1.030351901977319,1.045877422653877,1.034521134953779)
    ,@(1.02,1.027017579802031,1.036339756493796,1.027106533899833,1.042712989996007,1.023820303822167,1.033153148943728,1.039648222307023,1.030446896955262,1.045999800914294,1.034620344919551)
    ,@(1.02,1.026448817456673,1.035780531310927,1.026615253236135,1.042100000733194,1.023846096030143,1.032709930638099,1.039154838703062,1.03002190973766,1.045452368221459,1.03417649719309)
    ,@(1.02,1.024597212535539,1.033960483302448,1.025016108558192,1.040105381459722,1.02392621804881,1.031266183391793,1.037548228271204,1.028637609787624,1.043670307689587,1.032730699660094)
    ,@(1.02,1.023431048698816,1.032814588148412,1.024009121324592,1.03884988880944,1.023973741382756,1.030356243381403,1.036536060063308,1.027765185449951,1.042548027008646,1.031819467430517)
    ,@(1.02,1.024049446870191,1.0334222006188,1.024543093780495,1.039515586047209,1.023948818288664,1.03083883155319,1.037072825528395,1.028227872457806,1.043143146614358,1.032302740932862)
    ,@(1.02,1.026479827160454,1.035811019128138,1.026642037723908,1.042133418114489,1.023844704313072,1.032734098719934,1.039181740171068,1.030045083477312,1.045482214610346,1.034200699596373)
    ,@(1.02,1.029291935457606,1.038576723521798,1.029071367504201,1.045165632324112,1.02371141742819,1.034924208806534,1.041620551846369,1.032145194206914,1.048189016084252,1.036393919890506)
)

$colIndexes = @(2,3,4,5,6,9,10,11,12,13,14)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $startRow + $i
    $rowData = $data[$i]
    for ($j = 0; $j -lt $colIndexes.Count; $j++) {
        $ws.Cells.Item($rowNum, $colIndexes[$j]).Value = $rowData[$j]
    }
}

Write-Host "Updated $($data.Count) rows"
